# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 203, shifting the existing
# rows 203-292 down to 204-293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 203 (shifts 203..292 -> 204..293)
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record's data
$ws.Range("A203").Value() = 5
$ws.Range("B203").Value() = "Macroferia Regional de Talca"
$ws.Range("C203").Value() = "Maule"
$ws.Range("D203").Value() = 44468
$ws.Range("E203").Value() = 7
$ws.Range("F203").Value() = "Fruta"
$ws.Range("G203").Value() = 100101
$ws.Range("H203").Value() = "Berries"
$ws.Range("I203").Value() = 100112025
$ws.Range("J203").Value() = "Frutilla"
$ws.Range("K203").Value() = "Sin especificar"
$ws.Range("L203").Value() = "Primera"
$ws.Range("M203").Value() = 1000
$ws.Range("N203").Value() = 15000
$ws.Range("O203").Value() = 15000
$ws.Range("P203").Value() = 15000
$ws.Range("Q203").Value() = "$/bandeja 7 kilos"
$ws.Range("R203").Value() = "Provincia de Melipilla"
$ws.Range("S203").Value() = 2143
$ws.Range("T203").Value() = 7
